$d = $word.ActiveDocument

# 1. Remove the stray "_GoBack" bookmark that sits at the end of the
#    title paragraph ("... / Practicas de ciclos con while").
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2. Remove exercises 4-6 (numbered list items) plus the two blank
#    "ListParagraph" spacer paragraphs that followed them. These are
#    paragraphs 7 through 11 (1-based) at this point in the document:
#      7  - "...desde el ultimo caracter de la cadena al primero."
#      8  - "...unicamente las vocales a y u..."
#      9  - "Vuelva a hacer el ejercicio 5..."
#      10 - (blank)
#      11 - (blank)
$startPara = $d.Paragraphs.Item(7)
$endPara = $d.Paragraphs.Item(11)
$rangeToDelete = $d.Range($startPara.Range.Start, $endPara.Range.End)
$rangeToDelete.Delete()

# 3. The paragraph that used to hold the "Extra: ..." palindrome
#    exercise is now paragraph 7; empty its text but keep the
#    paragraph (and its ListParagraph/spacing/justify formatting).
$lastPara = $d.Paragraphs.Item(7)
$textRange = $lastPara.Range
[void]$textRange.MoveEnd(1, -1)
$textRange.Text = ""

# 4. Re-create the "_GoBack" bookmark, now collapsed inside that
#    emptied final paragraph.
$lastPara2 = $d.Paragraphs.Item(7)
$lastPara2.Range.Bookmarks.Add("_GoBack")
